$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-14 Thursday" "2025-08-15 Friday"

Replace-Text "93÷9=10, 3" "72÷4=18, 0"
Replace-Text "75÷8=9, 3" "15÷2=7, 1"
Replace-Text "49÷2=24, 1" "95÷6=15, 5"
Replace-Text "31÷2=15, 1" "12÷3=4, 0"
Replace-Text "63÷6=10, 3" "84÷5=16, 4"
Replace-Text "46÷9=5, 1" "40÷9=4, 4"
Replace-Text "92÷8=11, 4" "12÷6=2, 0"
Replace-Text "18÷8=2, 2" "13÷2=6, 1"
Replace-Text "43÷6=7, 1" "58÷5=11, 3"
Replace-Text "94÷4=23, 2" "44÷4=11, 0"
Replace-Text "58÷6=9, 4" "24÷7=3, 3"
Replace-Text "71÷8=8, 7" "62÷7=8, 6"
Replace-Text "27÷2=13, 1" "27÷9=3, 0"
Replace-Text "20÷8=2, 4" "30÷8=3, 6"
Replace-Text "11÷7=1, 4" "15÷6=2, 3"
Replace-Text "71÷4=17, 3" "59÷6=9, 5"
Replace-Text "81÷5=16, 1" "12÷7=1, 5"
Replace-Text "89÷3=29, 2" "84÷7=12, 0"
Replace-Text "29÷7=4, 1" "52÷8=6, 4"
Replace-Text "74÷6=12, 2" "28÷4=7, 0"
Replace-Text "90÷7=12, 6" "57÷5=11, 2"
Replace-Text "72÷9=8, 0" "67÷3=22, 1"
Replace-Text "36÷7=5, 1" "76÷7=10, 6"
Replace-Text "48÷2=24, 0" "31÷3=10, 1"
Replace-Text "37÷9=4, 1" "61÷7=8, 5"

Write-Output "Done"
